$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Samples")

# Update the flow channel naming format from hyphens to underscores.
$ws.Range("C3").Value = "FSC_H,SSC_H,FL1_H,FL1_H,FL3_H,FL1_A,FL4_H"

# Reflect the user's final navigation state: Samples tab active, C10 selected.
$ws.Activate() | Out-Null
$ws.Range("C10").Select() | Out-Null
